# Update the cryptocurrency price/volume table with freshly scraped values.
# Cells in column D that look like plain numbers would otherwise be
# auto-converted from text to a numeric type by Excel, so those are entered
# with a leading apostrophe to force them to stay as text (exactly like a
# user typing '22.13 into a cell does). Column D values that already contain
# more than one "." (e.g. 26.911.46) are never auto-converted, so they are
# assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.911.46"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").Value = "1.549.96"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").Value = "'206.57"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").Value = "'0.487"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").Value = "'22.13"
$ws.Range("E8").Value = "  +2.95%  "

$ws.Range("D9").Value = "'0.246"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").Value = "1.770.88"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").Value = "1.551.20"
$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("E14").Value = "  +0.65%  "

$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").Value = "26.896.99"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").Value = "'61.65"
$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").Value = "'217.31"
$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("E19").Value = "  +1.38%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("D23").Value = "'9.20"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("E24").Value = "  -0.65%  "

$ws.Range("E25").Value = "  +0.41%  "

$ws.Range("D26").Value = "'6.62"
$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("E28").Value = "  +0.78%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("D33").Value = "1.415.43"
$ws.Range("E33").Value = "  +3.15%  "

$ws.Range("E34").Value = "  +4.24%  "

$ws.Range("E35").Value = "  +2.51%  "

$ws.Range("D36").Value = "'0.970"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").Value = "'0.524"
$ws.Range("E39").Value = "  +1.10%  "

$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").Value = "'5.75"
$ws.Range("E41").Value = "  +4.79%  "

$ws.Range("D43").Value = "'2.32"
$ws.Range("E43").Value = "  +4.06%  "

$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  +1.59%  "

$ws.Range("D45").Value = "'64.25"
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("D46").Value = "'1.74"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Value = "1.684.45"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").Value = "'87.46"
$ws.Range("E48").Value = "  +1.42%  "

$ws.Range("E49").Value = "  +1.99%  "

$ws.Range("E50").Value = "  +3.98%  "

$ws.Range("D51").Value = "'0.0955"
$ws.Range("E51").Value = "  -0.27%  "
